# "added colors to rows" -- color-code specific DTR rows and fix a couple
# of incidental issues (FLOOR() called with an invalid 3rd arg, and a
# leave-day flag) that rode along with the same edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blue fill for the rows on 2015-02-10 (row 5) and 2015-02-20 (row 15)
$blue = 0xCCA329
$ws.Range("A5:J5").Interior.Color = $blue
$ws.Range("A15:J15").Interior.Color = $blue

# Red fill for the sick-leave row (2015-02-19, row 14) and flag it as a
# sick-leave day in the SICK LEAVE column.
$red = 0x5E5EDF
$ws.Range("A14:J14").Interior.Color = $red
$ws.Range("I14").Value = 1

# Tidy up the FLOOR(...) calls -- FLOOR only takes two arguments; the
# stray third argument was making these formulas error out.
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
